# Update countries & provincias Spain
# Applies the 30-Abr-2020 17:22 data refresh to the "Pais" sheet:
#  - updates the timestamp footer string
#  - refreshes case counts for several countries
#  - re-sorts a few country pairs whose totals crossed each other
#  - appends a new country row (Comoras)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Footer timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 30 de Abril de 2020 a las 17:22"

# --- Plain value refreshes (country/order unchanged) -------------------
$ws.Range("B14").Value = 80246
$ws.Range("C14").Value = 885
$ws.Range("E14").Value = 40573
$ws.Range("G14").Value = 30
$ws.Range("H14").Value = 5541

$ws.Range("E37").Value = 7518
$ws.Range("G37").Value = 12
$ws.Range("H37").Value = 705

$ws.Range("B44").Value = 7738
$ws.Range("C44").Value = 28
$ws.Range("E44").Value = 7499

$ws.Range("B65").Value = 2591
$ws.Range("C65").Value = 15
$ws.Range("D65").Value = 1374
$ws.Range("E65").Value = 1077
$ws.Range("F65").Value = 38
$ws.Range("G65").Value = 1
$ws.Range("H65").Value = 140

$ws.Range("B79").Value = 1506
$ws.Range("C79").Value = 59
$ws.Range("E79").Value = 1174
$ws.Range("F79").Value = 40
$ws.Range("G79").Value = 2
$ws.Range("H79").Value = 66

$ws.Range("B89").Value = 1089
$ws.Range("C89").Value = 12
$ws.Range("D89").Value = 642
$ws.Range("E89").Value = 445

$ws.Range("B94").Value = 850
$ws.Range("C94").Value = 7
$ws.Range("E94").Value = 687

$ws.Range("B128").Value = 315
$ws.Range("C128").Value = 2
$ws.Range("E128").Value = 36

# --- Rows 47/48: Australia and Republica Dominicana swap places --------
# Republica Dominicana moves up to row 47 with refreshed totals,
# Australia drops to row 48 keeping its previous totals.
$ws.Range("A47").Value = "Republica Dominicana"
$ws.Range("B47").Value = 6972
$ws.Range("C47").Value = 320
$ws.Range("D47").Value = 1301
$ws.Range("E47").Value = 5370
$ws.Range("F47").Value = 144
$ws.Range("G47").Value = 8
$ws.Range("H47").Value = 301

$ws.Range("A48").Value = "Australia"
$ws.Range("B48").Value = 6753
$ws.Range("C48").Value = 7
$ws.Range("D48").Value = 5715
$ws.Range("E48").Value = 947
$ws.Range("F48").Value = 34
$ws.Range("G48").Value = 2
$ws.Range("H48").Value = 91

# --- Rows 58/59: Argelia and Moldavia swap places -----------------------
# Argelia moves up to row 58 with refreshed totals,
# Moldavia drops to row 59 keeping its previous totals.
$ws.Range("A58").Value = "Argelia"
$ws.Range("B58").Value = 4006
$ws.Range("C58").Value = 158
$ws.Range("D58").Value = 1779
$ws.Range("E58").Value = 1777
$ws.Range("F58").Value = 22
$ws.Range("G58").Value = 6
$ws.Range("H58").Value = 450

$ws.Range("A59").Value = "Moldavia"
$ws.Range("B59").Value = 3897
$ws.Range("C59").Value = 126
$ws.Range("D59").Value = 1182
$ws.Range("E59").Value = 2599
$ws.Range("F59").Value = 237
$ws.Range("G59").Value = 5
$ws.Range("H59").Value = 116

# --- Rows 191/192: Namibia and San Vicente y las Granadinas swap places -
# Totals are identical for both countries, only the names trade places.
$ws.Range("A191").Value = "Namibia"
$ws.Range("A192").Value = "San Vicente y las Granadinas"

# --- New row 218: Comoras -----------------------------------------------
$ws.Range("A218").Value = "Comoras"
$ws.Range("B218").Value = 1
$ws.Range("C218").Value = 1
$ws.Range("D218").Value = 0
$ws.Range("E218").Value = 1
$ws.Range("F218").Value = 0
$ws.Range("G218").Value = 0
$ws.Range("H218").Value = 0
